$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.036.47'
$ws.Range("E2").Value = '  +0.62%  '

$ws.Range("D3").Value = '1.683.40'
$ws.Range("E3").Value = '  +0.73%  '

$ws.Range("E4").Value = '  -0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '215.98'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("E6").Value = '  -2.96%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.252'
$c.ClearFormats()
$ws.Range("E8").Value = '  -1.44%  '

$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '21.38'
$c.ClearFormats()
$ws.Range("E9").Value = '  +5.15%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0621'
$c.ClearFormats()
$ws.Range("E10").Value = '  +0.14%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0887'
$c.ClearFormats()

$ws.Range("D12").Value = '1.919.16'
$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").Value = '1.688.69'
$ws.Range("E13").Value = '  +0.76%  '

$ws.Range("E14").Value = '  +0.19%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.535'
$c.ClearFormats()
$ws.Range("E15").Value = '  +1.80%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '66.21'
$c.ClearFormats()
$ws.Range("E16").Value = '  +0.79%  '

$ws.Range("D17").Value = '27.050.41'
$ws.Range("E17").Value = '  +0.60%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '8.21'
$c.ClearFormats()
$ws.Range("E18").Value = '  +4.48%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '236.34'
$c.ClearFormats()
$ws.Range("E19").Value = '  +1.49%  '

$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("E22").Value = '  -0.16%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.25'
$c.ClearFormats()
$ws.Range("E23").Value = '  +0.42%  '

$ws.Range("E24").Value = '  -3.68%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '146.84'
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '7.24'
$c.ClearFormats()
$ws.Range("E26").Value = '  +1.30%  '

$ws.Range("E27").Value = '  +0.40%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.112'
$c.ClearFormats()
$ws.Range("E28").Value = '  -3.44%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("E31").Value = '  -0.39%  '

$ws.Range("E32").Value = '  -0.05%  '

$ws.Range("D33").Value = '1.514.39'
$ws.Range("E33").Value = '  +3.29%  '

$ws.Range("E34").Value = '  +0.22%  '

$ws.Range("E35").Value = '  +4.15%  '

$ws.Range("E36").Value = '  -0.54%  '

$ws.Range("E37").Value = '  +3.33%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.921'
$c.ClearFormats()
$ws.Range("E38").Value = '  +1.87%  '

$ws.Range("E39").Value = '  +2.83%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.04'
$c.ClearFormats()
$ws.Range("E40").Value = '  +7.43%  '

$ws.Range("E41").Value = '  -3.77%  '

$ws.Range("E42").Value = '  -0.09%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '68.39'
$c.ClearFormats()
$ws.Range("E43").Value = '  +3.88%  '

$ws.Range("E44").Value = '  -1.13%  '

$ws.Range("D45").Value = '1.822.30'
$ws.Range("E45").Value = '  +0.02%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.783'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.34%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '90.12'
$c.ClearFormats()
$ws.Range("E47").Value = '  -0.38%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.105'
$c.ClearFormats()
$ws.Range("E48").Value = '  +4.48%  '

$ws.Range("E49").Value = '  -0.86%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.86'
$c.ClearFormats()
$ws.Range("E50").Value = '  +3.56%  '

$ws.Range("E51").Value = '  +0.01%  '
